$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = '60.861.57'
$ws.Range("E2").Value = '  -3.27%  '
$ws.Range("D3").Value = '2.916.99'
$ws.Range("E3").Value = '  -3.85%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'584.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.52%  '
$ws.Range("D6").Value = "'144.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.10%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -2.64%  '
$ws.Range("D9").Value = '2.915.65'
$ws.Range("E9").Value = '  -3.66%  '
$ws.Range("D10").Value = "'6.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.18%  '
$ws.Range("E11").Value = '  -3.87%  '
$ws.Range("E12").Value = '  -4.08%  '
$ws.Range("E13").Value = '  -3.24%  '
$ws.Range("D14").Value = "'33.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.41%  '
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("D16").Value = '3.399.78'
$ws.Range("E16").Value = '  -3.87%  '
$ws.Range("D17").Value = '60.823.16'
$ws.Range("E17").Value = '  -3.31%  '
$ws.Range("E18").Value = '  -4.84%  '
$ws.Range("D19").Value = '2.918.19'
$ws.Range("E19").Value = '  -3.77%  '
$ws.Range("D20").Value = "'431.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = '  -4.56%  '
$ws.Range("E22").Value = '  -1.67%  '
$ws.Range("E23").Value = '  -4.45%  '
$ws.Range("D24").Value = "'80.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.19%  '
$ws.Range("D25").Value = "'10.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.40%  '
$ws.Range("E26").Value = '  -4.30%  '
$ws.Range("D27").Value = "'11.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.39%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("D30").Value = "'7.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.49%  '
$ws.Range("E31").Value = '  -3.08%  '
$ws.Range("D32").Value = "'2.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.40%  '
$ws.Range("E33").Value = '  -3.62%  '
$ws.Range("D34").Value = "'0.106"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.81%  '
$ws.Range("D35").Value = '0.0₃0870'
$ws.Range("E35").Value = '  +0.64%  '
$ws.Range("E36").Value = '  -2.60%  '
$ws.Range("E37").Value = '  -4.33%  '
$ws.Range("D38").Value = "'3.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.78%  '
$ws.Range("D39").Value = "'49.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.58%  '
$ws.Range("E40").Value = '  -2.11%  '
$ws.Range("E41").Value = '  -4.56%  '
$ws.Range("E42").Value = '  -4.28%  '
$ws.Range("E43").Value = '  -4.15%  '
$ws.Range("D44").Value = "'41.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.28%  '
$ws.Range("D45").Value = "'0.0348"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.90%  '
$ws.Range("D46").Value = "'375.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.97%  '
$ws.Range("D47").Value = '2.673.71'
$ws.Range("E47").Value = '  -2.21%  '
$ws.Range("D48").Value = "'132.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("D50").Value = "'24.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.73%  '
$ws.Range("E51").Value = '  -1.61%  '
